$d = $word.ActiveDocument

$replacements = @(
    @("54÷3=", "23÷2="),
    @("10÷4=", "64÷9="),
    @("18÷6=", "25÷9="),
    @("25÷3=", "88÷2="),
    @("14÷3=", "22÷3="),
    @("54÷4=", "40÷9="),
    @("16÷9=", "11÷6="),
    @("74÷2=", "47÷9="),
    @("64÷3=", "49÷2="),
    @("49÷8=", "14÷5="),
    @("50÷4=", "43÷8="),
    @("95÷6=", "28÷3="),
    @("31÷7=", "33÷6="),
    @("87÷9=", "95÷3="),
    @("41÷6=", "32÷8="),
    @("21÷3=", "12÷4="),
    @("89÷9=", "61÷5="),
    @("86÷8=", "90÷5="),
    @("99÷3=", "58÷7="),
    @("11÷5=", "93÷6="),
    @("33÷7=", "16÷2="),
    @("14÷6=", "34÷6="),
    @("44÷2=", "38÷6="),
    @("19÷6=", "72÷2="),
    @("22÷7=", "51÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
